$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D51").Value = "[github] public 레포지토리 private으로 변경하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/github-public-%EB%A0%88%ED%8F%AC%EC%A7%80%ED%86%A0%EB%A6%AC-private%EC%9C%BC%EB%A1%9C-%EB%B3%80%EA%B2%BD%ED%95%98%EA%B8%B0"

$ws.Range("D52").Value = "Relative Risk Regression"
